$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text values (prices/percentages) are written as literal strings,
# not auto-converted to numbers, by forcing a text number format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.202.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.425.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.420.69"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.871.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.968.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.426.56"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.92"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.44"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "631.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.540.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0943"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -9.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.25%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.94"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.75%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.44"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.373"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "148.10"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.27"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.22"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.45"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.63"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0518"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.592"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.53"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0235"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.07%  "
